$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the runs of a whole paragraph (identified by its trimmed
# text) with several new runs, described as an array of literal strings.
# Uses Range.InsertXML so that each array entry becomes its own <w:r>,
# matching the way Word splits a sentence into separate runs when it is
# edited in multiple steps.
# ---------------------------------------------------------------------------
function Set-ParagraphRuns($matchText, [string[]]$newRuns) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $matchText) {
            $r = $p.Range
            # Exclude the trailing paragraph mark from the replaced range.
            $textRange = $d.Range($r.Start, $r.End - 1)

            $runsXml = ""
            foreach ($piece in $newRuns) {
                $escaped = $piece.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
                if ($piece.Length -eq 0 -or $piece -ne $piece.Trim()) {
                    $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
                } else {
                    $runsXml += '<w:r><w:t>' + $escaped + '</w:t></w:r>'
                }
            }

            $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
                '</w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'

            $textRange.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# 1) "Names for placed orders." -> split into 4 runs.
[void](Set-ParagraphRuns "Names for placed orders. " @(
    "First or last name only ",
    "for placed orders",
    " (not both)",
    ". "
))

# 2) "Additional toppings..." -> split into 3 runs.
[void](Set-ParagraphRuns "Additional toppings, and taco/burrito shells will be available a la carte." @(
    "Additional toppings, and taco/burrito shells will be ",
    "complimentary",
    " a la carte."
))

# 3) Add a new bullet "Change/cash back" right after "Refills", at the same
#    outline level (ListParagraph / ilvl 1 / numId 1). Skip if it is already
#    there (keeps the script safe to re-run).
$alreadyThere = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Refills") {
        $nextText = $p.Next().Range.Text.TrimEnd([char]13)
        if ($nextText -eq "Change/cash back") {
            $alreadyThere = $true
        }
        break
    }
}

if (-not $alreadyThere) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq "Refills") {
            $p.Range.InsertParagraphAfter()
            break
        }
    }

    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq "Refills") {
            $newPara = $p.Next()
            $nr = $newPara.Range
            $target = $d.Range($nr.Start, $nr.End - 1)
            $target.Text = "Change/cash back"
            break
        }
    }
}
